$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell A1: "name" -> "command"
$ws.Range("A1").Value = "command"

# Widen column A to fit the longer susan command strings
$ws.Columns.Item(1).ColumnWidth = 60.15

# New row 7: automotive/susan smoothing
$ws.Range("A7").Value = "automotive/susan/susan input_large.pgm output_large.smoothing.pgm -s"
$ws.Range("B7").Value = 0.02
$ws.Range("C7").Value = 0.02
$ws.Range("D7").Value = 0

# New row 8: automotive/susan edges
$ws.Range("A8").Value = "automotive/susan/susan input_large.pgm output_large.edges.pgm -e"
$ws.Range("B8").Value = 0.01
$ws.Range("C8").Value = 0.01
$ws.Range("D8").Value = 0

# New row 9: automotive/susan corners
$ws.Range("A9").Value = "automotive/susan/susan input_large.pgm output_large.corners.pgm -c"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

# Match the alternating-row style already used on A3/A6 (copy formatting only)
$ws.Range("A8").Style = $ws.Range("A3").Style
$ws.Range("A9").Style = $ws.Range("A3").Style
$ws.Range("B9").Style = $ws.Range("A3").Style
$ws.Range("C9").Style = $ws.Range("A3").Style
$ws.Range("D9").Style = $ws.Range("A3").Style

$ws.Range("A23").Select()
